# ---------------------------------------------------------------
# Weekly CompStat refresh: new crime data collected (week of 6/16/2025-6/22/2025)
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and the week-covering date range ---
$ws.Range("A8").Value = "Volume 32   Number  25"
$ws.Range("C9").Value = "Report Covering the Week  6/16/2025  Through  6/22/2025"

# --- Template cells used to (re)apply the correct numeric/text style ---
# (style 13 = text, right-aligned placeholder used for "0" / "***.*";
#  style 14 = "#,##0" counts; style 15 = "#,##0.0;"-"#,##0.0" % change)
$styleTextCell = $ws.Range("C15")
$styleCountCell = $ws.Range("I15")
$stylePctCell = $ws.Range("L15")

# --- CompStat weekly crime table updates (rows 15-31) ---
$ws.Range("D15").Value = 1
$styleCountCell.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$stylePctCell.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$styleCountCell.Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$stylePctCell.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = 50
$ws.Range("N15").Value = -43.75
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 9.523809523809
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 113
$ws.Range("K16").Value = -15.04424778761
$ws.Range("L16").Value = -17.948717948717
$ws.Range("M16").Value = -8.571428571428
$ws.Range("N16").Value = -74.535809018567
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -70.588235294117
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = -50.980392156862
$ws.Range("I17").Value = 149
$ws.Range("J17").Value = 240
$ws.Range("K17").Value = -37.916666666666
$ws.Range("L17").Value = -18.579234972677
$ws.Range("M17").Value = 29.565217391304
$ws.Range("N17").Value = -49.662162162162
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "'0"
$styleTextCell.Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = "'***.*"
$styleTextCell.Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -18.181818181818
$ws.Range("I18").Value = 62
$ws.Range("K18").Value = 26.530612244898
$ws.Range("L18").Value = 14.814814814814
$ws.Range("M18").Value = 24
$ws.Range("N18").Value = -74.380165289256
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = -31.372549019607
$ws.Range("I19").Value = 192
$ws.Range("J19").Value = 252
$ws.Range("K19").Value = -23.809523809523
$ws.Range("L19").Value = -7.692307692307
$ws.Range("M19").Value = 88.235294117647
$ws.Range("N19").Value = 14.970059880239
$ws.Range("C20").Value = 3
$styleCountCell.Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 9
$ws.Range("I20").Value = 39
$ws.Range("J20").Value = 59
$ws.Range("K20").Value = -33.898305084745
$ws.Range("L20").Value = 18.181818181818
$ws.Range("M20").Value = 56
$ws.Range("N20").Value = -75.316455696202
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -46.153846153846
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 144
$ws.Range("H21").Value = -31.25
$ws.Range("I21").Value = 547
$ws.Range("J21").Value = 725
$ws.Range("K21").Value = -24.551724137931
$ws.Range("L21").Value = -9.586776859504
$ws.Range("M21").Value = 33.0900243309
$ws.Range("N21").Value = -56.963021243115
$ws.Range("G22").Value = "'0"
$styleTextCell.Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "'***.*"
$styleTextCell.Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("L22").Value = -15.384615384615
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 23.076923076923
$ws.Range("I23").Value = 81
$ws.Range("J23").Value = 84
$ws.Range("K23").Value = -3.571428571428
$ws.Range("L23").Value = 8
$ws.Range("M23").Value = 72.340425531914
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 21.428571428571
$ws.Range("F24").Value = 61
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -36.458333333333
$ws.Range("I24").Value = 377
$ws.Range("J24").Value = 469
$ws.Range("K24").Value = -19.616204690831
$ws.Range("L24").Value = -37.582781456953
$ws.Range("M24").Value = -11.084905660377
$ws.Range("D25").Value = 1
$styleCountCell.Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("E25").Value = 200
$stylePctCell.Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = 37.5
$ws.Range("I25").Value = 94
$ws.Range("J25").Value = 99
$ws.Range("K25").Value = -5.050505050505
$ws.Range("L25").Value = -66.187050359712
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -22.222222222222
$ws.Range("F26").Value = 56
$ws.Range("G26").Value = 78
$ws.Range("H26").Value = -28.205128205128
$ws.Range("I26").Value = 305
$ws.Range("J26").Value = 362
$ws.Range("K26").Value = -15.745856353591
$ws.Range("L26").Value = 24.489795918367
$ws.Range("M26").Value = 39.269406392694
$ws.Range("D27").Value = 1
$styleCountCell.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$stylePctCell.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 28.571428571428
$ws.Range("I28").Value = 32
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 6.666666666666
$ws.Range("L28").Value = 18.518518518518
$ws.Range("D29").Value = "'0"
$styleTextCell.Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = "'***.*"
$styleTextCell.Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 1
$ws.Range("M29").Value = -80
$ws.Range("N29").Value = -82.608695652173
$ws.Range("D30").Value = "'0"
$styleTextCell.Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$styleTextCell.Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 1
$ws.Range("M30").Value = -76.470588235294
$ws.Range("N30").Value = -82.608695652173
$ws.Range("D31").Value = "'0"
$styleTextCell.Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$styleTextCell.Copy()
$ws.Range("E31").PasteSpecial(-4122)
